$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# --- Insert a new row for "Manjaro 20.2 x86_64 / VMware" right after the openSUSE row (37) ---
$ws.Rows("38:38").Insert()
$ws.Range("H38").Clear()

$ws.Range("A38").Value = "Manjaro"
$ws.Range("B38").Value = 20.2
$ws.Range("C38").Value = "x86_64"
$ws.Range("E38").Value = "gcc 10.2.0"
$ws.Range("F38").Value = "VMware"
$ws.Range("G38").Value = 44215

# --- Insert a new row for "Manjaro 20.12 AArch64 / Raspberry Pi 4B 4GB" after the Raspbian
#     armv6l row (which, after the insert above, is now row 43) ---
$ws.Rows("44:44").Insert()
$ws.Range("H44").Clear()

$ws.Range("A44").Value = "Manjaro"
$ws.Range("B44").Value = 20.12
# B44 inherited a date number format from the row above when the row was inserted;
# pull the plain "General/centered" format from its neighbour instead.
$ws.Range("D44").Copy()
$ws.Range("B44").PasteSpecial(-4122)
$ws.Range("C44").Value = "AArch64"
$ws.Range("E44").Value = "gcc 10.2.0"
$ws.Range("F44").Value = "Raspberry Pi 4B 4GB"
$ws.Range("G44").Value = 44215

# --- Leave the window scrolled down with H17 selected, matching the author's screen state ---
$wb.Windows.Item(1).ScrollRow = 16
$wb.Windows.Item(1).ScrollColumn = 1
$ws.Range("H17").Select()
